$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.393.52"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "1.858.71"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.81"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3805"
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07306"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9296"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.74"
$ws.Range("E11").Value = "  +4.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07792"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "1.860.49"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.440"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.536"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.91"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008809"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "27.416.56"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.095"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.43"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.46"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.003"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.930"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08888"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.323"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.206"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.589"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7482"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.713"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.121"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  +3.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5527"
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.988"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05254"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.012"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.585"
$ws.Range("E42").Value = "  +4.03%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4869"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.660"
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.76"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.28"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06087"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9108"
$ws.Range("E51").Value = "  +2.37%  "
